$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.401.17'

$ws.Range("E2").Value = '  +2.52%  '

$ws.Range("D3").Value = '2.064.39'

$ws.Range("E3").Value = '  +3.91%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.80'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -0.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.614'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  +2.71%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.24'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  +7.01%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  +3.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.54'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  +0.05%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0763'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  +2.00%  '

$ws.Range("E12").Value = '  +3.31%  '

$ws.Range("D13").Value = '2.369.14'

$ws.Range("E13").Value = '  +3.92%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.55'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +1.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.02'
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.780'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +3.56%  '

$ws.Range("D18").Value = '2.047.94'

$ws.Range("E18").Value = '  +2.49%  '

$ws.Range("D19").Value = '37.599.14'

$ws.Range("E19").Value = '  +3.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.22'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +18.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.11'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +2.08%  '

$ws.Range("D22").Value = '0.0₃0818'

$ws.Range("E22").Value = '  +1.83%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '226.74'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +2.29%  '

$ws.Range("E24").Value = '  +0.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.46'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  +2.09%  '

$ws.Range("E26").Value = '  +0.99%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.36'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  +1.13%  '

$ws.Range("E28").Value = '  +14.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.89'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  +3.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.21'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  +2.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.127'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  -0.52%  '

$ws.Range("E32").Value = '  +2.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.50'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +2.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0624'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +2.91%  '

$ws.Range("E35").Value = '  +9.42%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.50'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  +6.08%  '

$ws.Range("E37").Value = '  +4.42%  '

$ws.Range("E38").Value = '  -0.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.79'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +0.48%  '

$ws.Range("E40").Value = '  +8.93%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0987'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +7.27%  '

$ws.Range("E42").Value = '  -1.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.48'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  +23.79%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '96.96'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +9.45%  '

$ws.Range("D45").Value = '1.473.40'

$ws.Range("E45").Value = '  +1.11%  '

$ws.Range("E46").Value = '  +6.19%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.96'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  +7.32%  '

$ws.Range("E49").Value = '  +3.66%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.29'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +7.70%  '

$ws.Range("E51").Value = '  +1.50%  '
